$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "for whom" (column E) values on existing rows
$ws.Range("E2").Value = "Мужчины"
$ws.Range("E14").Value = "Инопланетянин"

# New row 16 - another offer for the same Antigel AVS AVK-167 product
$ws.Range("A16").Value = "Авто и мото"
$ws.Range("B16").Value = "Автохимия"
$ws.Range("C16").Value = "Зимняя автохимия"
$ws.Range("D16").Value = "Антигель"
$ws.Range("E16").Value = "Инопланетянин212312"
$ws.Range("F16").Value = "Антигель AVS AVK-167, 340 мл 2 шт"
$ws.Range("G16").Value = 269
$ws.Range("H16").Value = "В наличии"
$ws.Range("I16").Value = 100
$ws.Range("J16").Value = "Предотвратите загустевание дизельного топлива на морозе с помощью антигеля AVS AVK-167. Он значительно облегчает пуск двигателя и снижает расход горючего, что продлевает срок службы механизма. Применяется для всех дизельных систем, в том числе Common rail и «насос-форсунка».`nСредство повышает работоспособность двигателя при температуре −31 °С при использовании летнего топлива и при −57 °С при использовании зимнего топлива.`nЗаливайте антигель в бак автомобиля перед заправкой. Одного флакона вам хватит на 60 –120 литров горючего.`nСостав: органический растворитель (> 30 %), депрессорная присадка (5 –15 %).`nОбъём: 520 мл."
$ws.Range("K16").Value = "https://cdn3.static1-sima-land.com/items/1593093/0/700-nw.jpg;https://cdn3.static1-sima-land.com/items/1593093/1/700-nw.jpg;https://cdn3.static1-sima-land.com/items/1593093/2/700-nw.jpg;https://cdn3.static1-sima-land.com/items/1593093/3/700-nw.jpg"
$ws.Range("L16").Value = "AVS"
$ws.Range("M16").Value = "https://www.sima-land.ru/1916276/"

# The multi-line description would otherwise force an auto row-height; reset it
# back to the default (no custom height) like every other data row.
$ws.Rows.Item(16).AutoFit()

# Update selection to mirror the saved view state
$ws.Range("E23").Select()
